$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Extend used range from column U to column AK for all rows (1-69) with blank cells.
$ws.Range("V1:AK69").Style = "Normal"

# 2. Update the "dct:modified" timestamp in B31.
$ws.Range("B31").Value = "2022-06-04T11:12:52+00:00"

# 3. Update header row 33: new iop/puv/sosa/owl columns are inserted after I33,
#    pushing the previous J/K/L header values (skos:editorialNote@en,
#    dct:creator(separator=","), dct:contributor(separator=",")) out to T/U/V.
$ws.Range("J33").Value = "iop:hasProperty"
$ws.Range("K33").Value = "iop:hasObjectOfInterest"
$ws.Range("L33").Value = "iop:hasMatrix"
$ws.Range("M33").Value = 'iop:hasContextObject(separator=",")'
$ws.Range("N33").Value = 'iop:hasConstraint(separator=",")'
$ws.Range("O33").Value = 'puv:statistic(separator=",")'
$ws.Range("P33").Value = 'puv:usesMethod(separator=",")'
$ws.Range("Q33").Value = 'sosa:madeBySensor(separator=",")'
$ws.Range("R33").Value = 'puv:uom(separator=",")'
$ws.Range("S33").Value = "owl:deprecated^^xsd:boolean"
$ws.Range("T33").Value = "skos:editorialNote@en"
$ws.Range("U33").Value = 'dct:creator(separator=",")'
$ws.Range("V33").Value = 'dct:contributor(separator=",")'

# 4. Row 34 (incentive-vars:1000 / IgG Subclasses):
#    C34 (qudt:unit "MFI") moves to R34 (now puv:uom column);
#    K34 (orcid, dct:creator) moves to U34 (now dct:creator column);
#    new iop:* columns J-Q get their values.
$ws.Range("C34").ClearContents()
$ws.Range("C34").Style = "Normal"
$ws.Range("J34").Value = "iop:hasProperty"
$ws.Range("K34").Value = "iop:hasObjectOfInterest"
$ws.Range("L34").Value = "iop:hasMatrix"
$ws.Range("M34").Value = 'iop:hasContextObject(separator=",")'
$ws.Range("N34").Value = 'iop:hasConstraint(separator=",")'
$ws.Range("O34").Value = 'puv:statistic(separator=",")'
$ws.Range("P34").Value = 'puv:usesMethod(separator=",")'
$ws.Range("Q34").Value = 'sosa:madeBySensor(separator=",")'
$ws.Range("R34").Value = "MFI"
$ws.Range("U34").Value = "https://orcid.org/0000-0003-3277-3107"

# 5. Rows 35-38 (IgG1-HA1 .. IgG4-HA1): qudt:unit ("MFI") moves C -> R,
#    dct:creator (orcid) moves K -> U.
foreach ($r in 35..38) {
    $ws.Range("C$r").ClearContents()
    $ws.Range("C$r").Style = "Normal"
    $ws.Range("K$r").ClearContents()
    $ws.Range("K$r").Style = "Normal"
    $ws.Range("R$r").Value = "MFI"
    $ws.Range("U$r").Value = "https://orcid.org/0000-0003-3277-3107"
}

# 6. Rows 39-43 (Surface Antigens-FC, CD3, Cytokines-FC, IFNg, Chemokines-FC):
#    qudt:unit ("MFI,%") moves C -> R.
foreach ($r in 39..43) {
    $ws.Range("C$r").ClearContents()
    $ws.Range("C$r").Style = "Normal"
    $ws.Range("R$r").Value = "MFI,%"
}
